$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# A leading apostrophe forces Excel to store numeric-looking values as text
# (preserving formatting such as trailing zeros), matching the source data.

$ws.Range("D2").Value = '31.601.71'
$ws.Range("E2").Value = '  +5.82%  '
$ws.Range("D3").Value = '1.710.30'
$ws.Range("E3").Value = '  +4.40%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '''221.96'
$ws.Range("E5").Value = '  +3.02%  '
$ws.Range("D6").Value = '''0.535'
$ws.Range("E6").Value = '  +3.14%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").Value = '''30.00'
$ws.Range("E8").Value = '  +4.28%  '
$ws.Range("D9").Value = '''45.37'
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("E10").Value = '  +3.65%  '
$ws.Range("E11").Value = '  +5.73%  '
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").Value = '1.954.79'
$ws.Range("E13").Value = '  +4.39%  '
$ws.Range("D14").Value = '1.710.95'
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").Value = '''10.27'
$ws.Range("E15").Value = '  +8.62%  '
$ws.Range("E16").Value = '  +3.63%  '
$ws.Range("E17").Value = '  +8.09%  '
$ws.Range("D18").Value = '31.530.51'
$ws.Range("E18").Value = '  +5.56%  '
$ws.Range("D19").Value = '''67.26'
$ws.Range("E19").Value = '  +4.17%  '
$ws.Range("D20").Value = '''250.71'
$ws.Range("E20").Value = '  +4.26%  '
$ws.Range("D21").Value = '0.0₃0725'
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").Value = '''0.997'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").Value = '''10.20'
$ws.Range("E23").Value = '  +2.90%  '
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("D26").Value = '''159.40'
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("E28").Value = '  +3.06%  '
$ws.Range("E29").Value = '  +2.85%  '
$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  +11.49%  '
$ws.Range("D32").Value = '''0.0504'
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("E33").Value = '  +3.80%  '
$ws.Range("D34").Value = '''3.41'
$ws.Range("E34").Value = '  +6.59%  '
$ws.Range("D35").Value = '1.517.35'
$ws.Range("E35").Value = '  +6.54%  '
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("D38").Value = '''83.60'
$ws.Range("E38").Value = '  +8.89%  '
$ws.Range("D39").Value = '''0.610'
$ws.Range("E39").Value = '  +8.73%  '
$ws.Range("E40").Value = '  +4.20%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").Value = '''2.31'
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("D44").Value = '''0.857'
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("D45").Value = '''0.0504'
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").Value = '''1.04'
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").Value = '''52.44'
$ws.Range("E48").Value = '  +7.37%  '
$ws.Range("D49").Value = '''5.57'
$ws.Range("E49").Value = '  +3.49%  '
$ws.Range("D50").Value = '1.843.43'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("E51").Value = '  +10.13%  '
